$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (TUE): Collaborative Development -> Human Computer Interaction, time shifts later
$ws.Range("B4").Value = "9:00-11:00"
$ws.Range("C4").Value = "5CS022"
$ws.Range("D4").Value = "Human Computer Interaction"
$ws.Range("G4").Value = "Mr. Pravash Karki"
$ws.Range("H4").Value = "LT-02 Telford"
$ws.Range("J4").Value = "L5CG(1+2+3+4)"

# Row 5 (TUE): Human Computer Interaction -> Distributed and Cloud Systems Programming
$ws.Range("B5").Value = "12:00-14:00"
$ws.Range("C5").Value = "5CS020"
$ws.Range("D5").Value = "Distributed and Cloud Systems Programming"
$ws.Range("G5").Value = "Mr. Sumanta Silwal"
$ws.Range("H5").Value = "LT-01 Wulfruna"

# Row 6: TUE -> WED, Distributed... -> Collaborative Development
$ws.Range("A6").Value = "WED"
$ws.Range("B6").Value = "9:30-11:30"
$ws.Range("C6").Value = "5CS024"
$ws.Range("D6").Value = "Collaborative Development"
$ws.Range("G6").Value = "Mr. Udaya Kandel"

# Row 7: Human Computer Interaction -> Distributed and Cloud Systems Programming, Lecture -> Tutorial
$ws.Range("B7").Value = "12:30-14:30"
$ws.Range("C7").Value = "5CS020"
$ws.Range("D7").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F7").Value = "Tutorial"
$ws.Range("G7").Value = "Mr. Prabin Sapkota"
$ws.Range("H7").Value = "TR-03 Westbromwich"
$ws.Range("J7").Value = "L5CG2"

# Row 8: WED -> THU, Collaborative Development -> Human Computer Interaction, Lecture -> Tutorial
$ws.Range("A8").Value = "THU"
$ws.Range("B8").Value = "9:00-11:00"
$ws.Range("C8").Value = "5CS022"
$ws.Range("D8").Value = "Human Computer Interaction"
$ws.Range("F8").Value = "Tutorial"
$ws.Range("G8").Value = "Mr. Dipesh Shrestha"
$ws.Range("H8").Value = "TR-02 Stafford"
$ws.Range("J8").Value = "L5CG2"

# Row 9: WED -> FRI, Distributed... -> Collaborative Development
$ws.Range("A9").Value = "FRI"
$ws.Range("B9").Value = "10:00-12:00"
$ws.Range("C9").Value = "5CS024"
$ws.Range("D9").Value = "Collaborative Development"
$ws.Range("G9").Value = "Mr. Anmol Adhikari"
$ws.Range("H9").Value = "TR-09  Chandragiri"
$ws.Range("I9").Value = "HCK"

# Row 10: THU -> FRI, Tutorial -> Workshop, hours 2 -> 2.5
$ws.Range("A10").Value = "FRI"
$ws.Range("B10").Value = "13:30-16:00"
$ws.Range("E10").Value = 2.5
$ws.Range("F10").Value = "Workshop"
$ws.Range("H10").Value = "SR-03 Wolves"

# Remove rows 11-13 (schedule now only has 9 class entries, dimension shrinks to A1:L10)
$ws.Rows("11:13").Delete()
